$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 439, all hold the "Förändrad" date value that was
# bumped by one day (date serial 45188 -> 45189).
$ws.Range("C2:C439").Value = 45189
